$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.893517000000001
$ws.Range("H2").Value = 26.680551
$ws.Range("I2").Value = 0.9082944842335181
$ws.Range("J2").Value = 0.9082944842335181
$ws.Range("M2").Value = 0.9705896666666667
$ws.Range("N2").Value = 2.911769
$ws.Range("O2").Value = 0.02073452941466921
$ws.Range("P2").Value = 0.02073452941466921
$ws.Range("Q2").Value = 8.631955700524335
$ws.Range("R2").Value = 77.68760130471901
$ws.Range("S2").Value = 0.01883305870052168
$ws.Range("T2").Value = 0.01883305870052168
$ws.Range("G3").Value = 8.893517000000001
$ws.Range("H3").Value = 26.680551
$ws.Range("I3").Value = 0.9082944842335181
$ws.Range("J3").Value = 0.9082944842335181
$ws.Range("O3").Value = 0.5628689972673966
$ws.Range("P3").Value = 0.5628689972673966
$ws.Range("Q3").Value = 234.32700846219
$ws.Range("R3").Value = 2108.94307615971
$ws.Range("S3").Value = 0.5112508055640275
$ws.Range("T3").Value = 0.5112508055640275
$ws.Range("G4").Value = 8.893517000000001
$ws.Range("H4").Value = 26.680551
$ws.Range("I4").Value = 0.9082944842335181
$ws.Range("J4").Value = 0.9082944842335181
$ws.Range("M4").Value = 19.49164633333333
$ws.Range("N4").Value = 58.47493899999999
$ws.Range("O4").Value = 0.4163964733179342
$ws.Range("P4").Value = 0.4163964733179341
$ws.Range("Q4").Value = 173.3492880234877
$ws.Range("R4").Value = 1560.143592211389
$ws.Range("S4").Value = 0.3782106199689689
$ws.Range("T4").Value = 0.3782106199689689
$ws.Range("I5").Value = 0.04237443292342908
$ws.Range("J5").Value = 0.04237443292342909
$ws.Range("M5").Value = 0.9705896666666667
$ws.Range("N5").Value = 2.911769
$ws.Range("O5").Value = 0.02073452941466921
$ws.Range("P5").Value = 0.02073452941466921
$ws.Range("Q5").Value = 0.4027044468276667
$ws.Range("R5").Value = 3.624340021449
$ws.Range("S5").Value = 0.0008786139258807678
$ws.Range("T5").Value = 0.000878613925880768
$ws.Range("I6").Value = 0.04237443292342908
$ws.Range("J6").Value = 0.04237443292342909
$ws.Range("O6").Value = 0.5628689972673966
$ws.Range("P6").Value = 0.5628689972673966
$ws.Range("S6").Value = 0.02385125456938508
$ws.Range("T6").Value = 0.02385125456938509
$ws.Range("I7").Value = 0.04237443292342908
$ws.Range("J7").Value = 0.04237443292342909
$ws.Range("M7").Value = 19.49164633333333
$ws.Range("N7").Value = 58.47493899999999
$ws.Range("O7").Value = 0.4163964733179342
$ws.Range("P7").Value = 0.4163964733179341
$ws.Range("Q7").Value = 8.087220505224332
$ws.Range("R7").Value = 72.78498454701899
$ws.Range("S7").Value = 0.01764456442816323
$ws.Range("T7").Value = 0.01764456442816323
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.4830226666666667
$ws.Range("H8").Value = 1.449068
$ws.Range("I8").Value = 0.04933108284305281
$ws.Range("J8").Value = 0.04933108284305281
$ws.Range("M8").Value = 0.9705896666666667
$ws.Range("N8").Value = 2.911769
$ws.Range("O8").Value = 0.02073452941466921
$ws.Range("P8").Value = 0.02073452941466921
$ws.Range("Q8").Value = 0.4688168090324444
$ws.Range("R8").Value = 4.219351281292
$ws.Range("S8").Value = 0.001022856788266762
$ws.Range("T8").Value = 0.001022856788266762
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.4830226666666667
$ws.Range("H9").Value = 1.449068
$ws.Range("I9").Value = 0.04933108284305281
$ws.Range("J9").Value = 0.04933108284305281
$ws.Range("O9").Value = 0.5628689972673966
$ws.Range("P9").Value = 0.5628689972673966
$ws.Range("Q9").Value = 12.72671503292
$ws.Range("R9").Value = 114.54043529628
$ws.Range("S9").Value = 0.02776693713398401
$ws.Range("T9").Value = 0.02776693713398401
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.4830226666666667
$ws.Range("H10").Value = 1.449068
$ws.Range("I10").Value = 0.04933108284305281
$ws.Range("J10").Value = 0.04933108284305281
$ws.Range("M10").Value = 19.49164633333333
$ws.Range("N10").Value = 58.47493899999999
$ws.Range("O10").Value = 0.4163964733179342
$ws.Range("P10").Value = 0.4163964733179341
$ws.Range("Q10").Value = 9.414906989650222
$ws.Range("R10").Value = 84.73416290685199
$ws.Range("S10").Value = 0.02054128892080204
$ws.Range("T10").Value = 0.02054128892080204
